# Insert a new data row at row 350 (shifting old rows 350-429 down to 351-430)
# and populate the new row 350 with the updated record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 350; this shifts rows 350..429 down to 351..430
$ws.Rows("350").Insert()

# Populate the newly inserted row 350
$ws.Range("A350").Value2 = 3
$ws.Range("B350").Value2 = "Femacal de La Calera"
$ws.Range("C350").Value2 = "Coquimbo"
$ws.Range("D350").Value2 = 45275
$ws.Range("E350").Value2 = 5
$ws.Range("F350").Value2 = "Fruta"
$ws.Range("G350").Value2 = 100101
$ws.Range("H350").Value2 = "Berries"
$ws.Range("I350").Value2 = 100101001
$ws.Range("J350").Value2 = "Arándano (blue)"
$ws.Range("K350").Value2 = "Sin especificar"
$ws.Range("L350").Value2 = "Primera"
$ws.Range("M350").Value2 = 40
$ws.Range("N350").Value2 = 6000
$ws.Range("O350").Value2 = 6000
$ws.Range("P350").Value2 = 6000
$ws.Range("Q350").Value2 = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R350").Value2 = "Provincia de Curicó"
$ws.Range("S350").Value2 = 4000
$ws.Range("T350").Value2 = 1.5
